$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-05-11 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-12 Monday", 2) | Out-Null

# Update the division-problem answers in the table, cell by cell
# (addressed by row/column so identical strings that appear as both
#  an old value in one cell and a new value in another never collide)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "57÷7=8, 1"
$t.Cell(1, 2).Range.Text = "74÷8=9, 2"
$t.Cell(1, 3).Range.Text = "94÷6=15, 4"
$t.Cell(1, 4).Range.Text = "53÷4=13, 1"
$t.Cell(1, 5).Range.Text = "53÷2=26, 1"

$t.Cell(5, 1).Range.Text = "81÷6=13, 3"
$t.Cell(5, 2).Range.Text = "45÷3=15, 0"
$t.Cell(5, 3).Range.Text = "72÷4=18, 0"
$t.Cell(5, 4).Range.Text = "70÷4=17, 2"
$t.Cell(5, 5).Range.Text = "89÷2=44, 1"

$t.Cell(9, 1).Range.Text = "28÷8=3, 4"
$t.Cell(9, 2).Range.Text = "73÷4=18, 1"
$t.Cell(9, 3).Range.Text = "76÷8=9, 4"
$t.Cell(9, 4).Range.Text = "30÷3=10, 0"
$t.Cell(9, 5).Range.Text = "66÷9=7, 3"

$t.Cell(13, 1).Range.Text = "90÷6=15, 0"
$t.Cell(13, 2).Range.Text = "50÷6=8, 2"
$t.Cell(13, 3).Range.Text = "24÷6=4, 0"
$t.Cell(13, 4).Range.Text = "64÷5=12, 4"
$t.Cell(13, 5).Range.Text = "97÷2=48, 1"

$t.Cell(17, 1).Range.Text = "64÷2=32, 0"
$t.Cell(17, 2).Range.Text = "21÷3=7, 0"
$t.Cell(17, 3).Range.Text = "64÷6=10, 4"
$t.Cell(17, 4).Range.Text = "14÷4=3, 2"
$t.Cell(17, 5).Range.Text = "52÷9=5, 7"
